$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 21:05"

$ws.Range("B4").Value = 1715799
$ws.Range("C4").Value = 9573
$ws.Range("D4").Value = 469064
$ws.Range("E4").Value = 1146519
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 411
$ws.Range("H4").Value = 100216

$ws.Range("B10").Value = 182722
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 65879
$ws.Range("E10").Value = 88313
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 98
$ws.Range("H10").Value = 28530

$ws.Range("B11").Value = 181236
$ws.Range("C11").Value = 447
$ws.Range("D11").Value = 162000
$ws.Range("E11").Value = 10760
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 48
$ws.Range("H11").Value = 8476

$ws.Range("B13").Value = 150772
$ws.Range("C13").Value = 5822
$ws.Range("D13").Value = 64277
$ws.Range("E13").Value = 82146
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 177
$ws.Range("H13").Value = 4349

$ws.Range("B114").Value = 956
$ws.Range("C114").Value = 5
$ws.Range("D114").Value = 634
$ws.Range("E114").Value = 312
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 10

$ws.Range("A124").Value = "Nicaragua"
$ws.Range("B124").Value = 759
$ws.Range("C124").Value = 480
$ws.Range("D124").Value = 370
$ws.Range("E124").Value = 354
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 18
$ws.Range("H124").Value = 35

$ws.Range("A125").Value = "Sierra Leona"
$ws.Range("B125").Value = 754
$ws.Range("C125").Value = 19
$ws.Range("D125").Value = 297
$ws.Range("E125").Value = 413
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 44

$ws.Range("A126").Value = "Georgia"
$ws.Range("B126").Value = 732
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 537
$ws.Range("E126").Value = 183
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 12

$ws.Range("A127").Value = "Jordania"
$ws.Range("B127").Value = 718
$ws.Range("C127").Value = 7
$ws.Range("D127").Value = 586
$ws.Range("E127").Value = 123
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 9

$ws.Range("A128").Value = "Crucero"
$ws.Range("B128").Value = 712
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 651
$ws.Range("E128").Value = 48
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 13

$ws.Range("A129").Value = "Etiopia"
$ws.Range("B129").Value = 701
$ws.Range("C129").Value = 46
$ws.Range("D129").Value = 167
$ws.Range("E129").Value = 528
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 6

$ws.Range("A130").Value = "Republica del Chad"
$ws.Range("B130").Value = 700
$ws.Range("C130").Value = 13
$ws.Range("D130").Value = 303
$ws.Range("E130").Value = 335
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 62

$ws.Range("A131").Value = "Republica de Africa Central"
$ws.Range("B131").Value = 671
$ws.Range("C131").Value = 19
$ws.Range("D131").Value = 22
$ws.Range("E131").Value = 648
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 1

$ws.Range("A132").Value = "San Marino"
$ws.Range("B132").Value = 666
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 275
$ws.Range("E132").Value = 349
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 42

$ws.Range("A133").Value = "Malta"
$ws.Range("B133").Value = 611
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 485
$ws.Range("E133").Value = 120
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 6

$ws.Range("A134").Value = "Madagascar"
$ws.Range("B134").Value = 586
$ws.Range("C134").Value = 44
$ws.Range("D134").Value = 147
$ws.Range("E134").Value = 437
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 2

$ws.Range("A135").Value = "Jamaica"
$ws.Range("B135").Value = 556
$ws.Range("C135").Value = 4
$ws.Range("D135").Value = 238
$ws.Range("E135").Value = 309
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 9

$ws.Range("A136").Value = "Tanzania"
$ws.Range("B136").Value = 509
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 183
$ws.Range("E136").Value = 305
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 21

$ws.Range("A137").Value = "Congo"
$ws.Range("B137").Value = 487
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 147
$ws.Range("E137").Value = 324
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 16

$ws.Range("A138").Value = "Reunion"
$ws.Range("B138").Value = 459
$ws.Range("C138").Value = 3
$ws.Range("D138").Value = 411
$ws.Range("E138").Value = 47
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 1

$ws.Range("A139").Value = "Taiwan"
$ws.Range("B139").Value = 441
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 416
$ws.Range("E139").Value = 18
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 7

$ws.Range("A140").Value = "Estado de Palestina"
$ws.Range("B140").Value = 426
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 365
$ws.Range("E140").Value = 58
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 3

$ws.Range("A141").Value = "Cabo Verde"
$ws.Range("B141").Value = 390
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 155
$ws.Range("E141").Value = 231
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 4

$ws.Range("A142").Value = "Togo"
$ws.Range("B142").Value = 386
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 161
$ws.Range("E142").Value = 212
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 13

$ws.Range("A143").Value = "Guayana Francesa"
$ws.Range("B143").Value = 353
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 146
$ws.Range("E143").Value = 206
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 1

$ws.Range("A144").Value = "Ruanda"
$ws.Range("B144").Value = 339
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 244
$ws.Range("E144").Value = 95
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

$ws.Range("A145").Value = "Isla de Man"
$ws.Range("B145").Value = 336
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 306
$ws.Range("E145").Value = 6
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 24

$ws.Range("A146").Value = "Mauricio"
$ws.Range("B146").Value = 334
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 322
$ws.Range("E146").Value = 2
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 10

$ws.Range("A147").Value = "Vietnam"
$ws.Range("B147").Value = 327
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 272
$ws.Range("E147").Value = 55
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0

$ws.Range("A148").Value = "Montenegro"
$ws.Range("B148").Value = 324
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 315
$ws.Range("E148").Value = 0
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 9

$ws.Range("A149").Value = "Santo Tome y Principe"
$ws.Range("B149").Value = 299
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 4
$ws.Range("E149").Value = 284
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 11

$ws.Range("A150").Value = "Mauritania"
$ws.Range("B150").Value = 268
$ws.Range("C150").Value = 6
$ws.Range("D150").Value = 15
$ws.Range("E150").Value = 240
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 4
$ws.Range("H150").Value = 13

$ws.Range("A151").Value = "Liberia"
$ws.Range("B151").Value = 266
$ws.Range("C151").Value = 1
$ws.Range("D151").Value = 144
$ws.Range("E151").Value = 96
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 26

$ws.Range("B190").Value = 25
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 18
$ws.Range("E190").Value = 6
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 1
